# Update the "Overall Demand" rows (2030, 2040, 2050) for the
# Iron & steel (C), Chemicals (D), and Non-metallic minerals (E) columns
# with the summed per-country demand values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2030 - row 2
$ws.Range("C2").Value = 0.005275666120273866
$ws.Range("D2").Value = 0.008851734721773526
$ws.Range("E2").Value = 0.01430471668413269

# 2040 - row 12
$ws.Range("C12").Value = 0.02218391830577439
$ws.Range("D12").Value = 0.01011012915445557
$ws.Range("E12").Value = 0.02168455665835297

# 2050 - row 22
$ws.Range("C22").Value = 0.05806349380617117
$ws.Range("D22").Value = 0.01132223274949354
$ws.Range("E22").Value = 0.03488305184526833
